# testReport2.xlsx — "import budget name from different column"
#
# The "Aufwände gesamt" sheet lists time-tracking records. Column E
# ("Gruppe") already carries the budget name (Budget1 / Budget2) for each
# row, but column F ("Subgruppe") - which the importer actually reads the
# budget name from - is empty for the imported rows. Populate F4:F11 with
# the same budget-name value already present in E4:E11 for that row, so the
# importer (which now reads the budget name from the different column) has
# data to read.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aufwände gesamt")
$ws.Activate()

for ($r = 4; $r -le 11; $r++) {
    $groupValue = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 6).Value2 = $groupValue
}

# Touch row 17 (just below the data) so it materializes as an explicit row
# in the sheet, matching the row the user's cursor had moved past.
$ws.Rows.Item(17).OutlineLevel = 0

# Leave the selection on the newly-filled column F range.
$ws.Range("F4:F11").Select()
